$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# A leading apostrophe forces text entry (matching the source data which
# stores these numbers as text, e.g. "36.482.94"), then the style is reset
# to Normal so no stray quote-prefix formatting remains on the cell.
$ws.Range("D2").Value = "'36.482.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.98%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.019.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.28%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'232.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -10.03%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.599"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.24%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.00%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'54.83"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.59%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -2.56%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'57.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.23%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0748"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.83%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.60%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'2.315.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.38%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'14.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.41%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'20.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -5.41%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.760"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -4.92%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -2.66%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.018.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.81%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'36.659.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.08%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'67.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -5.32%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.0₃0797"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.58%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +5.41%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'220.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -5.65%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.04%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +1.17%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -7.73%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'162.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.70%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.133"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +4.72%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'8.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -4.10%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.23%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'18.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.49%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.116"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.35%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -5.11%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.0599"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -6.83%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +4.11%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'4.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -4.05%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.17%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.95%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -6.73%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +3.46%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0966"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +4.16%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -4.15%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.455.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.11%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'FTXToken"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'4.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +38.85%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'VeChain"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'0.0204"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -3.15%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'90.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.93%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -6.99%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'15.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.52%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -1.90%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -1.73%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'6.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.33%  "
$ws.Range("E51").Style = "Normal"
